# Update job titles and "Last Bonus Allocation Percent" (column L) values
# on the active sheet, per the engineering-focused titles rename.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - Paige Duty: Staff Engineer -> Staff SRE; L3 95 -> 115
$ws.Range("C3").Value = "Staff SRE"
$ws.Range("L3").Value = 115

# Row 4 - Lee Latency: Senior Software Engineer -> Senior Software Developer; L4 new value 90
$ws.Range("C4").Value = "Senior Software Developer"
$ws.Range("L4").Value = 90

# Row 5 - Mona Torr: Senior Software Engineer -> Senior SRE; L5 90 -> 85
$ws.Range("C5").Value = "Senior SRE"
$ws.Range("L5").Value = 85

# Row 6 - Robin Rollback: Software Engineer -> Software Developer; L6 new value 115
$ws.Range("C6").Value = "Software Developer"
$ws.Range("L6").Value = 115

# Row 7 - Kenny Canary: Software Engineer -> Software Developer; L7 115 -> 85
$ws.Range("C7").Value = "Software Developer"
$ws.Range("L7").Value = 85

# Row 8 - Tracey Loggins: Senior Software Engineer -> Senior SRE; L8 95 -> 90
$ws.Range("C8").Value = "Senior SRE"
$ws.Range("L8").Value = 90

# Row 9 - Sue Q. Ell: Senior Software Engineer -> Senior Software Developer; L9 unchanged (105)
$ws.Range("C9").Value = "Senior Software Developer"

# Row 10 - Jason Blob: Software Engineer -> Software Developer; L10 cleared (cell removed)
$ws.Range("C10").Value = "Software Developer"
$ws.Range("L10").ClearContents()

# Row 11 - Al Ert: Staff Engineer -> Staff SRE; L11 110 -> 95
$ws.Range("C11").Value = "Staff SRE"
$ws.Range("L11").Value = 95

# Row 12 - Addie Min: Senior Software Engineer -> Senior Software Developer; L12 115 -> 105
$ws.Range("C12").Value = "Senior Software Developer"
$ws.Range("L12").Value = 105

# Row 13 - Tim Out: Software Engineer -> Software Developer; L13 100 -> 110
$ws.Range("C13").Value = "Software Developer"
$ws.Range("L13").Value = 110

# Row 14 - Barbie Que: Senior Software Engineer -> Senior SRE; L14 115 -> 105
$ws.Range("C14").Value = "Senior SRE"
$ws.Range("L14").Value = 105
